$wb = $excel.ActiveWorkbook

$headers = @("frequency", "text", "lemma", "pos", "eng-tag", "dependency", "afinn sentiment", "mcdonals sentiment", "token id")

foreach ($ws in $wb.Worksheets) {
    for ($i = 0; $i -lt $headers.Length; $i++) {
        $col = 2 + $i  # Column B is index 2
        $ws.Cells.Item(1, $col).Value = $headers[$i]
    }
}
